$d = $word.ActiveDocument

# Find the placeholder run(s) "#panitia" ... "pejabat2#" (split across two
# runs, with the "_GoBack" bookmark sitting between them) and figure out
# which paragraph contains it.
$rng = $d.Content
$found = $rng.Find.Execute("#panitia*pejabat*#", $true, $false, $true, $false,
                            $false, $true, 1, $false, "", 0)

$target = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $rng.Start -and $p.Range.End -gt $rng.Start) {
        $target = $i
        break
    }
}

$p = $d.Paragraphs.Item($target)

# Delete the paragraph's content together with its paragraph mark so it
# merges into the following (currently empty) paragraph. The following
# paragraph's pPr - tabs/spacing/ind, with no centered jc - becomes the
# pPr of the merged paragraph.
$delRange = $d.Range($p.Range.Start, $p.Range.End)
$delRange.Delete()

# The merged, now-empty paragraph sits at the same index.
$merged = $d.Paragraphs.Item($target)

# Move the left indent from 4680 twips (234 pt) to 6840 twips (342 pt).
$merged.Format.LeftIndent = 342

# Re-insert the "_GoBack" bookmark at the start of the (empty) paragraph.
$d.Bookmarks.Add("_GoBack", $merged.Range)

# Insert the combined placeholder text, formatted to match the original
# runs (Tahoma, black, 11pt / 22 half-points, cyan highlight).
$ins = $merged.Range
$ins.InsertAfter("#panitiapejabat#")
$ins.Font.NameAscii = "Tahoma"
$ins.Font.Name = "Tahoma"
$ins.Font.Color = 0
$ins.Font.Size = 11
$ins.HighlightColorIndex = 3
